# "Add files via upload" - populate sheet "17" (sheet7.xml) with four new
# names, pulled in from the shared-strings pool, change the workbook's
# default font from Calibri to Arial, and make sheet "17" the active tab
# (selecting cell B7 on it).

$wb = $excel.ActiveWorkbook

# --- sheet "17" is the 7th tab (sheetId 7 / r:id rId7 -> sheet7.xml) ---
$ws7 = $wb.Worksheets.Item(7)

# Write the new names. The order below matters for shared-string allocation:
# A1, A2 get the first two new strings; A4 grabs the third (Noam Raanan)
# before A3 grabs the fourth (Lior Tsalovich), matching the final layout
# A1=Yuval Koskas, A2=Max Gutnik, A3=Lior Tsalovich, A4=Noam Raanan.
$ws7.Range("A1").Value = "Yuval Koskas "
$ws7.Range("A2").Value = "Max Gutnik"
$ws7.Range("A4").Value = "Noam Raanan"
$ws7.Range("A3").Value = "Lior Tsalovich"

# Workbook-wide default font: Calibri -> Arial
$wb.Styles.Item("Normal").Font.Name = "Arial"

# Make sheet "17" the active / selected tab, with B7 selected (was D7).
$ws7.Activate()
$ws7.Range("B7").Select()
